# Apply the "updated errors and gedcom" edit:
#  - Backlog sheet: flip several status cells in column E between
#    "coding" and "done", and fill in the previously-blank status
#    cells for rows 31, 40 and 41.
#  - Sprint2 sheet: flip status cells in column D from "coding" to
#    "done", and fill in the previously-blank G2/H2 estimate cells.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Backlog sheet
# ---------------------------------------------------------------
$backlog = $wb.Worksheets.Item("Backlog")

# coding -> done
$backlog.Range("E12").Value = "done"
$backlog.Range("E21").Value = "done"

# done -> coding
$doneToCoding = @(22,23,24,25,28,29,30,31,32,33,34,35,36,37,38,39,40,41)
foreach ($r in $doneToCoding) {
    $backlog.Range("E$r").Value = "coding"
}

# ---------------------------------------------------------------
# Sprint2 sheet
# ---------------------------------------------------------------
$sprint2 = $wb.Worksheets.Item("Sprint2")

# coding -> done
$codingToDone = @(2,3,4,9,10,11)
foreach ($r in $codingToDone) {
    $sprint2.Range("D$r").Value = "done"
}

# newly filled in estimate cells
$sprint2.Range("G2").Value = 50
$sprint2.Range("H2").Value = 50

# ---------------------------------------------------------------
# Selections left behind by the edit session.
# Backlog stays the active/tab-selected sheet, so select on Sprint2
# first and finish on Backlog.
# ---------------------------------------------------------------
$sprint2.Range("I16").Select() | Out-Null
$backlog.Range("F33").Select() | Out-Null

Write-Host "edit applied"
